$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter 4 into Z8. This ripples through the shared formulas across the rest
# of row 8 (AA8, AD8, AG8, AJ8, AM8, AP8, AS8, AV8, AY8, AZ8, BA8), which all
# recalculate automatically.
$ws.Range("Z8").Value = 4

# Re-merge these header cells (row 4) so they move to the end of the
# worksheet's merged-cell list, mirroring how Excel appends freshly
# (re)merged ranges to the end of <mergeCells>.
$headerMerges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($mergeRef in $headerMerges) {
    $ws.Range($mergeRef).UnMerge()
    $ws.Range($mergeRef).Merge()
}

# Move the active selection (bottom-right frozen pane) to E10.
$ws.Range("E10").Select()
